$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.149.18"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.601.69"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'580.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "'191.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "3.597.36"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "'0.666"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "'55.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("E13").Value = "  +5.93%  "
$ws.Range("D14").Value = "'9.72"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "4.187.47"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'19.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").Value = "3.609.17"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "70.196.30"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'483.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'19.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.35%  "
$ws.Range("E24").Value = "  -6.80%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "'95.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.36%  "
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").Value = "'11.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'9.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'32.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "'12.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'66.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "'588.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.62%  "
$ws.Range("D36").Value = "'39.08"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "0.0₃0805"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'3.33"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +22.43%  "
$ws.Range("D41").Value = "'3.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("E42").Value = "  -5.95%  "
$ws.Range("D43").Value = "3.236.37"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "'2.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.83%  "
$ws.Range("D45").Value = "'3.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'0.0450"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("D48").Value = "'3.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'3.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.66%  "
